# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E30) listed periods in descending order
# (2103 .. 2001). The update re-lists them in ascending order (2001 .. 2103),
# and the "Valor Mora" (F column) amounts follow their respective period -
# which, since only the F16/F30 rows actually differed (25396 vs 33125),
# shows up as those two values swapping places.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
if (-not $ws) { $ws = $wb.ActiveSheet }

# New ascending period sequence for rows 16..30.
$periods = @("2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# Valor Mora amounts follow the period, not the row, so the two rows whose
# amount differed from the common 33125 swap:  2001 -> 33125, 2103 -> 25396.
$ws.Cells.Item(16, 6).Value = 33125
$ws.Cells.Item(30, 6).Value = 25396
